$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.616.93'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.66%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.280.18'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.96%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("B5").Value = 'Solana'
$ws.Range("C5").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '187.96'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.05%  '
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '585.91'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.43%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.603'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.95%  '
$ws.Range("E9").Value = '  +0.22%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.64'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.19%  '
$ws.Range("E11").Value = '  -1.66%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.855.85'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.04%  '
$ws.Range("E13").Value = '  +1.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.65'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.58%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '68.619.04'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.65%  '
$ws.Range("E16").Value = '  -0.83%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.282.67'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.72%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.75'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.48'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '416.63'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +5.18%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.58'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.54%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.56'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.511'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.41%  '
$ws.Range("E25").Value = '  -0.36%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.191'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.49'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.01'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.85%  '
$ws.Range("E29").Value = '  -0.91%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.80'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.42%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.52'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.58%  '
$ws.Range("E32").Value = '  -0.73%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.90'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.58%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '164.21'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.22%  '
$ws.Range("E35").Value = '  -2.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.91'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.52%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.58'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.43%  '
$ws.Range("B38").Value = 'Filecoin'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.53'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.85%  '
$ws.Range("B39").Value = 'Mantle'
$ws.Range("C39").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.794'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -3.11%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.41'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.87%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.682.22'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.58%  '
$ws.Range("B42").Value = 'Hedera'
$ws.Range("C42").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0682'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.96%  '
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.45'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.71%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.57'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.34%  '
$ws.Range("B45").Value = 'Bittensor'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '338.14'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.36%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '24.80'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0276'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.44%  '
$ws.Range("B48").Value = 'Arweave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '32.07'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +2.72%  '
$ws.Range("B49").Value = 'ONDO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.996'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.17%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.20'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.49%  '
$ws.Range("E51").Value = '  -1.31%  '
